# edit.ps1
# Applies the "updated docx title page" edit:
#   1. Remove the stray "_GoBack" bookmark from the empty paragraph above
#      the title line.
#   2. Change the morning service time from 08.45 to 09.30.
#   3. Change the location from "Opstandingskerk - Wilgenbeemd 2" to
#      "Het Noorderlicht - Scharmbarg 37", and re-create the "_GoBack"
#      bookmark right after the new location text (mirroring where Word
#      leaves it after the last edit made to the document).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. It will be re-added below
#    after the location text is rewritten, since that is where Word
#    leaves this auto-maintained bookmark after the most recent edit.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Morning service time: 08.45 -> 09.30
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Aanvang morgendienst 08.45 uur", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Aanvang morgendienst 09.30 uur", 2)

# ---------------------------------------------------------------------
# 3a. Location: Opstandingskerk - Wilgenbeemd 2  ->  Het Noorderlicht - Scharmbarg 37
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Locatie: Opstandingskerk", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Locatie: Het Noorderlicht", 2)
$d.Content.Find.Execute("Wilgenbeemd 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Scharmbarg 37", 2)

# ---------------------------------------------------------------------
# 3b. Re-insert the "_GoBack" bookmark immediately after the new
#     location text (still inside that paragraph, before its paragraph
#     mark).
#
#     Adding a zero-length bookmark directly at the position right
#     before a paragraph mark is unreliable, so a single placeholder
#     character is inserted, the bookmark is created around it, and the
#     placeholder is then removed again - leaving a clean, collapsed
#     bookmark in the correct spot without altering the visible text.
# ---------------------------------------------------------------------
$locatiePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Scharmbarg 37*") {
        $locatiePara = $p
    }
}
$pos = $locatiePara.Range.End - 1

$placeholderChar = [char]1
$insertPoint = $d.Range($pos, $pos)
$insertPoint.InsertAfter($placeholderChar)

$placeholderRange = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$d.Bookmarks("_GoBack").Range.Text = ""
